# Generate Report for Handoff
# Adds two new source files (14432ec2-975e-438b-aa83-997f69c30a47 and
# 996f4713-836c-43b1-a933-fc8151987c3a) to the localization status report,
# inserting them (alphabetically, by GUID) between the existing
# 0b85b24d... row and the f532fd09... row on every sheet, and pushing the
# ".localization-config" row down accordingly.

$wb = $excel.ActiveWorkbook

$mdBase   = "https://github.com/OpenLocalizationTest/oltest/blob/36fd24e5b1ae63b1938796a6fc7b5ab5250f0030/e2e/"
$cfgUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/36fd24e5b1ae63b1938796a6fc7b5ab5250f0030/.localization-config"
$zhcnBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/518e1056381509209ba280ab5762b7ec1cbbb8d4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/"
$dedeBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b7ebbc5a50379b88deb3ed64df41c7aa1c6fa0b1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/"

$file14432 = "14432ec2-975e-438b-aa83-997f69c30a47.md"
$file996f  = "996f4713-836c-43b1-a933-fc8151987c3a.md"
$file0b85  = "0b85b24d-0558-4027-81c5-745b20487a57.md"
$filef532  = "f532fd09-7c1e-4e0d-a364-2c7f2da79398.md"
$cfgName   = ".localization-config"

$xlf14432zh = "14432ec2-975e-438b-aa83-997f69c30a47.f16d3f3c41a0779f1c591dd07979a72c63a34c0b.zh-cn.xlf"
$xlf996fzh  = "996f4713-836c-43b1-a933-fc8151987c3a.75cd08d881f6e052c06d7e16d9ee041a18cc8bb7.zh-cn.xlf"
$xlf0b85zh  = "0b85b24d-0558-4027-81c5-745b20487a57.420ad34371b22ac3378c66b77d6e6f000de8cf4a.zh-cn.xlf"
$xlff532zh  = "f532fd09-7c1e-4e0d-a364-2c7f2da79398.86bb0105902eae0d79f9ccb90c27e96e78c2a96a.zh-cn.xlf"

$xlf14432de = "14432ec2-975e-438b-aa83-997f69c30a47.f16d3f3c41a0779f1c591dd07979a72c63a34c0b.de-de.xlf"
$xlf996fde  = "996f4713-836c-43b1-a933-fc8151987c3a.75cd08d881f6e052c06d7e16d9ee041a18cc8bb7.de-de.xlf"
$xlf0b85de  = "0b85b24d-0558-4027-81c5-745b20487a57.420ad34371b22ac3378c66b77d6e6f000de8cf4a.de-de.xlf"
$xlff532de  = "f532fd09-7c1e-4e0d-a364-2c7f2da79398.86bb0105902eae0d79f9ccb90c27e96e78c2a96a.de-de.xlf"

$readyForHandoff = "Ready for handoff"
$notLocalized    = "Not to be localized"
$include         = "Include"
$ignored         = "Ignored"
$epoch           = "0001-01-01 00:00:00"

$dt14432zh = "2016-03-10 03:35:57"
$dt996fzh  = "2016-03-10 03:35:57"
$dt0b85zh  = "2016-03-10 03:35:20"
$dtf532zh  = "2016-03-10 03:35:20"

$dt14432de = "2016-03-10 03:36:05"
$dt996fde  = "2016-03-10 03:36:05"
$dt0b85de  = "2016-03-10 03:35:40"
$dtf532de  = "2016-03-10 03:35:40"

# -----------------------------------------------------------------
# Sheet 1: "Overview"
# -----------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Rows.Item(3).Insert()
$ws1.Rows.Item(3).Insert()

$ws1.Range("A3").Value = $file14432
$ws1.Range("B3").Value = $readyForHandoff
$ws1.Range("C3").Value = $readyForHandoff

$ws1.Range("A4").Value = $file996f
$ws1.Range("B4").Value = $readyForHandoff
$ws1.Range("C4").Value = $readyForHandoff

$ws1.Range("A3").Style = "Hyperlink"
$ws1.Range("A4").Style = "Hyperlink"

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), $mdBase + $file0b85, "", "", $file0b85)
$ws1.Hyperlinks.Add($ws1.Range("A3"), $mdBase + $file14432, "", "", $file14432)
$ws1.Hyperlinks.Add($ws1.Range("A4"), $mdBase + $file996f, "", "", $file996f)
$ws1.Hyperlinks.Add($ws1.Range("A5"), $mdBase + $filef532, "", "", $filef532)
$ws1.Hyperlinks.Add($ws1.Range("A6"), $cfgUrl, "", "", $cfgName)

# -----------------------------------------------------------------
# Sheet 2: "zh-cn"
# -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Rows.Item(3).Insert()
$ws2.Rows.Item(3).Insert()

$ws2.Range("A3").Value = $file14432
$ws2.Range("B3").Value = $readyForHandoff
$ws2.Range("C3").Value = $xlf14432zh
$ws2.Range("D3").Value = $dt14432zh
$ws2.Range("G3").Value = $epoch
$ws2.Range("H3").Value = $include

$ws2.Range("A4").Value = $file996f
$ws2.Range("B4").Value = $readyForHandoff
$ws2.Range("C4").Value = $xlf996fzh
$ws2.Range("D4").Value = $dt996fzh
$ws2.Range("G4").Value = $epoch
$ws2.Range("H4").Value = $include

$ws2.Range("A3").Style = "Hyperlink"
$ws2.Range("A4").Style = "Hyperlink"
$ws2.Range("C3").Style = "Hyperlink"
$ws2.Range("C4").Style = "Hyperlink"
$ws2.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $mdBase + $file0b85, "", "", $file0b85)
$ws2.Hyperlinks.Add($ws2.Range("C2"), $zhcnBase + $xlf0b85zh, "", "", $xlf0b85zh)
$ws2.Hyperlinks.Add($ws2.Range("A3"), $mdBase + $file14432, "", "", $file14432)
$ws2.Hyperlinks.Add($ws2.Range("C3"), $zhcnBase + $xlf14432zh, "", "", $xlf14432zh)
$ws2.Hyperlinks.Add($ws2.Range("A4"), $mdBase + $file996f, "", "", $file996f)
$ws2.Hyperlinks.Add($ws2.Range("C4"), $zhcnBase + $xlf996fzh, "", "", $xlf996fzh)
$ws2.Hyperlinks.Add($ws2.Range("A5"), $mdBase + $filef532, "", "", $filef532)
$ws2.Hyperlinks.Add($ws2.Range("C5"), $zhcnBase + $xlff532zh, "", "", $xlff532zh)
$ws2.Hyperlinks.Add($ws2.Range("A6"), $cfgUrl, "", "", $cfgName)

# -----------------------------------------------------------------
# Sheet 3: "de-de"
# -----------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Rows.Item(3).Insert()
$ws3.Rows.Item(3).Insert()

$ws3.Range("A3").Value = $file14432
$ws3.Range("B3").Value = $readyForHandoff
$ws3.Range("C3").Value = $xlf14432de
$ws3.Range("D3").Value = $dt14432de
$ws3.Range("G3").Value = $epoch
$ws3.Range("H3").Value = $include

$ws3.Range("A4").Value = $file996f
$ws3.Range("B4").Value = $readyForHandoff
$ws3.Range("C4").Value = $xlf996fde
$ws3.Range("D4").Value = $dt996fde
$ws3.Range("G4").Value = $epoch
$ws3.Range("H4").Value = $include

$ws3.Range("A3").Style = "Hyperlink"
$ws3.Range("A4").Style = "Hyperlink"
$ws3.Range("C3").Style = "Hyperlink"
$ws3.Range("C4").Style = "Hyperlink"
$ws3.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $mdBase + $file0b85, "", "", $file0b85)
$ws3.Hyperlinks.Add($ws3.Range("C2"), $dedeBase + $xlf0b85de, "", "", $xlf0b85de)
$ws3.Hyperlinks.Add($ws3.Range("A3"), $mdBase + $file14432, "", "", $file14432)
$ws3.Hyperlinks.Add($ws3.Range("C3"), $dedeBase + $xlf14432de, "", "", $xlf14432de)
$ws3.Hyperlinks.Add($ws3.Range("A4"), $mdBase + $file996f, "", "", $file996f)
$ws3.Hyperlinks.Add($ws3.Range("C4"), $dedeBase + $xlf996fde, "", "", $xlf996fde)
$ws3.Hyperlinks.Add($ws3.Range("A5"), $mdBase + $filef532, "", "", $filef532)
$ws3.Hyperlinks.Add($ws3.Range("C5"), $dedeBase + $xlff532de, "", "", $xlff532de)
$ws3.Hyperlinks.Add($ws3.Range("A6"), $cfgUrl, "", "", $cfgName)

Write-Host "Report regenerated for handoff."
